$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range('D2').Value = '29.546.76'
$ws.Range('E2').Value = '  +2.30%  '

# Row 3
$ws.Range('D3').Value = '1.990.85'
$ws.Range('E3').Value = '  +6.02%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').Value = "'325.42"
$ws.Range('E5').Value = '  +0.17%  '

# Row 6
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  -0.10%  '

# Row 7
$ws.Range('D7').Value = "'0.4681"
$ws.Range('E7').Value = '  +1.72%  '

# Row 8
$ws.Range('D8').Value = "'0.3946"
$ws.Range('E8').Value = '  +1.73%  '

# Row 9
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'46.47"
$ws.Range('E9').Value = '  -0.26%  '

# Row 10
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.07955"
$ws.Range('E10').Value = '  +1.31%  '

# Row 11
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = "'1.002"
$ws.Range('E11').Value = '  +1.67%  '

# Row 12
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = "'23.04"
$ws.Range('E12').Value = '  +5.93%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '2.002.90'
$ws.Range('E13').Value = '  +4.89%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = "'7.266"
$ws.Range('E14').Value = '  +3.92%  '

# Row 15
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = "'5.880"
$ws.Range('E15').Value = '  +4.20%  '

# Row 16
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = "'0.07150"
$ws.Range('E16').Value = '  +2.66%  '

# Row 17
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = "'88.75"
$ws.Range('E17').Value = '  +0.95%  '

# Row 18
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value = "'1.003"
$ws.Range('E18').Value = '  +0.05%  '

# Row 19
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = "'0.000009953"
$ws.Range('E19').Value = '  -0.38%  '

# Row 20
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = "'17.43"
$ws.Range('E20').Value = '  +2.65%  '

# Row 21
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = "'1.002"
$ws.Range('E21').Value = '  +0.03%  '

# Row 22
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '29.647.88'
$ws.Range('E22').Value = '  +2.66%  '

# Row 23
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = "'5.543"
$ws.Range('E23').Value = '  +6.02%  '

# Row 24
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = "'11.28"
$ws.Range('E24').Value = '  +3.04%  '

# Row 25
$ws.Range('D25').Value = "'2.106"
$ws.Range('E25').Value = '  +0.89%  '

# Row 26
$ws.Range('D26').Value = "'157.92"
$ws.Range('E26').Value = '  +1.14%  '

# Row 27
$ws.Range('D27').Value = "'19.64"
$ws.Range('E27').Value = '  +1.64%  '

# Row 28
$ws.Range('D28').Value = "'5.970"
$ws.Range('E28').Value = '  -0.83%  '

# Row 29
$ws.Range('D29').Value = "'120.41"
$ws.Range('E29').Value = '  +2.62%  '

# Row 30
$ws.Range('D30').Value = "'1.966"
$ws.Range('E30').Value = '  +2.10%  '

# Row 31
$ws.Range('D31').Value = "'0.09453"
$ws.Range('E31').Value = '  +1.05%  '

# Row 32
$ws.Range('D32').Value = "'0.9030"
$ws.Range('E32').Value = '  +0.35%  '

# Row 33
$ws.Range('D33').Value = "'5.259"
$ws.Range('E33').Value = '  +0.15%  '

# Row 34
$ws.Range('E34').Value = '  +2.49%  '

# Row 35
$ws.Range('D35').Value = "'3.180"
$ws.Range('E35').Value = '  -2.30%  '

# Row 36
$ws.Range('D36').Value = "'0.05845"
$ws.Range('E36').Value = '  +1.60%  '

# Row 37
$ws.Range('D37').Value = "'1.178"
$ws.Range('E37').Value = '  -0.29%  '

# Row 38
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'0.02124"
$ws.Range('E38').Value = '  +2.72%  '

# Row 39
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = "'0.000003304"
$ws.Range('E39').Value = '  +102.34%  '

# Row 40
$ws.Range('D40').Value = "'7.889"
$ws.Range('E40').Value = '  +3.08%  '

# Row 41
$ws.Range('D41').Value = "'0.5761"
$ws.Range('E41').Value = '  +2.18%  '

# Row 42
$ws.Range('D42').Value = "'0.1829"
$ws.Range('E42').Value = '  +3.79%  '

# Row 43
$ws.Range('D43').Value = "'9.824"
$ws.Range('E43').Value = '  +1.82%  '

# Row 44
$ws.Range('D44').Value = "'12.04"
$ws.Range('E44').Value = '  +1.73%  '

# Row 45
$ws.Range('D45').Value = "'0.5372"
$ws.Range('E45').Value = '  +0.69%  '

# Row 46
$ws.Range('D46').Value = "'2.698"
$ws.Range('E46').Value = '  +6.49%  '

# Row 47
$ws.Range('D47').Value = "'2.184"
$ws.Range('E47').Value = '  -3.83%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = "'1.870"
$ws.Range('E48').Value = '  +1.64%  '

# Row 49
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.06952"
$ws.Range('E49').Value = '  -1.30%  '

# Row 50
$ws.Range('D50').Value = "'114.61"
$ws.Range('E50').Value = '  +1.85%  '

# Row 51
$ws.Range('D51').Value = "'0.3084"
$ws.Range('E51').Value = '  +8.29%  '
